$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Cadastrado" (B) and "Sem Cadastro" (C) values for years 2021-2025
$ws.Range("B3").Value = 909716.14
$ws.Range("C3").Value = 46448.4

$ws.Range("B4").Value = 1749677.29
$ws.Range("C4").Value = 22690.39

$ws.Range("B5").Value = 2823826.46
$ws.Range("C5").Value = 21376.22

$ws.Range("B6").Value = 4420458.1
$ws.Range("C6").Value = 26899.54

$ws.Range("B7").Value = 2259434.43
$ws.Range("C7").Value = 13169.03
